# Add data for 2022-11-22
# Updates the YTD violent-crime counts (one additional day of incidents)
# across the citywide roll-up sheets and the affected neighborhood sheets.

$wb = $excel.ActiveWorkbook

# Citywide Totals: Aggravated Battery (row3), Robbery (row6), Total (row7)
$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("F2").Value = 87
$ws.Range("D3").Value = 130
$ws.Range("E3").Value = 139
$ws.Range("F3").Value = 131
$ws.Range("I3").Value = 186
$ws.Range("C6").Value = 460
$ws.Range("D6").Value = 394
$ws.Range("F6").Value = 499
$ws.Range("G6").Value = 425
$ws.Range("I6").Value = 483
$ws.Range("C7").Value = 611
$ws.Range("D7").Value = 620
$ws.Range("E7").Value = 665
$ws.Range("F7").Value = 725
$ws.Range("G7").Value = 651
$ws.Range("I7").Value = 807

# Little Italy, UIC
$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("D5").Value = 10
$ws.Range("D6").Value = 13

# Humboldt Park
$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("G4").Value = 5
$ws.Range("G5").Value = 10

# South Chicago
$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("D4").Value = 4
$ws.Range("D5").Value = 5

# South Shore
$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("F2").Value = 2
$ws.Range("F5").Value = 12

# By Neighborhood roll-up (Total row per neighborhood, plus grand Total row 98)
$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("I19").Value = 22
$ws.Range("G41").Value = 10
$ws.Range("E47").Value = 17
$ws.Range("D50").Value = 13
$ws.Range("I52").Value = 7
$ws.Range("C53").Value = 53
$ws.Range("F68").Value = 4
$ws.Range("D75").Value = 5
$ws.Range("F77").Value = 19
$ws.Range("I77").Value = 47
$ws.Range("D80").Value = 5
$ws.Range("F82").Value = 12
$ws.Range("I85").Value = 5
$ws.Range("C98").Value = 611
$ws.Range("D98").Value = 620
$ws.Range("E98").Value = 665
$ws.Range("F98").Value = 725
$ws.Range("G98").Value = 651
$ws.Range("I98").Value = 807

# Loop
$ws = $wb.Worksheets.Item('Loop')
$ws.Range("C6").Value = 35
$ws.Range("C7").Value = 53

# Roseland
$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("F3").Value = 7
$ws.Range("I3").Value = 9
$ws.Range("F7").Value = 19
$ws.Range("I7").Value = 47

# Chatham
$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("I5").Value = 14
$ws.Range("I6").Value = 22

# United Center
$ws = $wb.Worksheets.Item('United Center')
$ws.Range("I4").Value = 2
$ws.Range("I5").Value = 5

# Logan Square
$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("I5").Value = 5
$ws.Range("I6").Value = 7

# Lake View
$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("E3").Value = 3
$ws.Range("E6").Value = 17

# Riverdale
$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("D3").Value = 3
$ws.Range("D5").Value = 5

# Bridgeport
$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("F5").Value = 3

# O'Hare (apostrophe escaped by doubling inside the single-quoted string)
$ws = $wb.Worksheets.Item('O''Hare')
$ws.Range("F6").Value = 4
